# Analyze swap part 1
# Fills in benchmark "number of swaps" data on the "Geral" sheet for the
# 10-element dataset (rows 3-11), corrects several totals further down the
# sheet (rows 14-33), and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Geral")

function Set-Row($rowNum, $values) {
    # $values is an ordered list of 7 numbers/$null for columns B..H
    $range = $ws.Range("B$rowNum`:H$rowNum")
    $range.NumberFormat = "0"
    $arr = New-Object 'object[,]' 1,7
    for ($i = 0; $i -lt 7; $i++) {
        if ($null -ne $values[$i]) {
            $arr[0,$i] = $values[$i]
        } else {
            $arr[0,$i] = ""
        }
    }
    $range.Value = $arr
    # Blank out any cell that should stay empty (assigning "" would create a
    # text cell, so clear those specifically instead).
    for ($i = 0; $i -lt 7; $i++) {
        if ($null -eq $values[$i]) {
            $ws.Cells.Item($rowNum, 2 + $i).ClearContents()
        }
    }
}

# ---------------------------------------------------------------------------
# Rows 3-11: "Aleatorio" (random) dataset, n = 10 .. 1000000 (columns B-H),
# only the first four columns (10, 100, 1000, 10000) have known results so
# far; larger sizes (100000, 500000, 1000000) are left blank for now.
# ---------------------------------------------------------------------------
Set-Row 3  @(16, 2371, 249127, 24901016, $null, $null, $null)
Set-Row 4  @(16, 2371, 249127, 24901016, $null, $null, $null)
Set-Row 5  @(12, 159,  2408,   31390,    $null, $null, $null)
Set-Row 6  @(12, 188,  2674,   34066,    $null, $null, $null)
Set-Row 7  @(16, 2371, 249127, 24901016, $null, $null, $null)
Set-Row 8  @(9,  314,  7190,   212680,   $null, $null, $null)
Set-Row 9  @(8,  93,   994,    9981,     $null, $null, $null)
Set-Row 10 @(34, 672,  9976,   133616,   $null, $null, $null)
Set-Row 11 @(29, 581,  9096,   124299,   $null, $null, $null)

# Stray formatting note left on an otherwise-empty row.
$ws.Range("B12").Font.Underline = $true

# ---------------------------------------------------------------------------
# Rows 14-22: "Crescente" (already sorted) dataset.
# ---------------------------------------------------------------------------
$ws.Range("G14").Value = 0
$ws.Range("G15").Value = 0

Set-Row 16 @(10, 100, 1000, 10000, 100000, 500000, 1000000)

$ws.Range("B17").Font.Underline = $true
$ws.Range("G17").Value = 262143
$ws.Range("H17").Value = 524287

$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 0

Set-Row 21 @(34, 672, 9976, 133616, 1668928, 9475712, 19951424)

$ws.Range("G22").Value = 9355690
$ws.Range("H22").Value = 19794250

# ---------------------------------------------------------------------------
# Rows 25-33: "Decrescente" (reverse sorted) dataset.
# ---------------------------------------------------------------------------
$ws.Range("H25").Value = 1783293664
$ws.Range("G26").Value = 445698416

Set-Row 27 @(9, 99, 999, 9999, 99999, 500000, 999999)

$ws.Range("G28").Value = 512142
$ws.Range("H28").Value = 1024286

$ws.Range("G29").Value = 445698416
$ws.Range("H29").Value = 1783293664

$ws.Range("G30").Value = 499071604
$ws.Range("H30").Value = 1962492188

$ws.Range("G31").Value = 250000
$ws.Range("H31").Value = 500000

Set-Row 32 @(34, 672, 9976, 133616, 1668928, 9475712, 19951424)

$ws.Range("F33").Value = 1497466
$ws.Range("G33").Value = 8668486
$ws.Range("H33").Value = 18333446

# ---------------------------------------------------------------------------
# Update the sheet selection/scroll position.
# ---------------------------------------------------------------------------
$ws.Range("E7").Select()
